$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio2")

# Update Tabella14 values (row 4: Game Lv 3)
$ws.Range("B4").Value = 2600
$ws.Range("C4").Value = 1700

# Update Tabella135 row 16: move the 200 Exp value into the Exp opt (C) column
$ws.Range("B16").Value = $null
$ws.Range("C16").Value = 200

# Update the selection shown on the sheet
$ws.Range("F3:F8").Select()

$wb.Save()
